$wb = $excel.ActiveWorkbook
$plate1 = $wb.Worksheets.Item("Plate1")
$plate1.Copy([System.Reflection.Missing]::Value, $plate1)
$newSheet = $wb.Worksheets.Item(3)
$newSheet.Name = "Plate2"
